{"js": "// The document ends with an empty paragraph immediately before the\n// section break. Fill it with the \"Version management...\" text (it stays\n// the same paragraph - we are not inserting a new one).\nconst body = context.document.body;\nconst lastParagraph = body.paragraphs.getLast();\n\nlastParagraph.insertText(\n  \"Version management allows us to manage projects that are developed by multiple people at the same time. It can manage the modification history of a file or project. View historical versions. Back up or restore previous versions.\",\n  \"Start\"\n);\n\nawait context.sync();\n", "ps1": "# Add the \"Version management...\" paragraph text to the trailing empty\n# paragraph that currently sits right before the section break (the last\n# paragraph in the document body).\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n$lastParagraph.Range.Text = \"Version management allows us to manage projects that are developed by multiple people at the same time. It can manage the modification history of a file or project. View historical versions. Back up or restore previous versions.\"\n"}
